$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round sumsq (C), statistic (E), and p.value (F) columns to 3
# significant figures for data rows 2-19 (see commit message:
# 'update significant digits in output tables').

$ws.Range("C2").Value = 0.141
$ws.Range("E2").Value = 45.2
$ws.Range("F2").Value = 0.0000000405

$ws.Range("C3").Value = 0.0343
$ws.Range("E3").Value = 11.0
$ws.Range("F3").Value = 0.00195

$ws.Range("C4").Value = 0.00614
$ws.Range("E4").Value = 1.96
$ws.Range("F4").Value = 0.169

$ws.Range("C5").Value = 2.16
$ws.Range("E5").Value = 11.4
$ws.Range("F5").Value = 0.0016

$ws.Range("C6").Value = 3.82
$ws.Range("E6").Value = 20.2
$ws.Range("F6").Value = 0.0000554

$ws.Range("C7").Value = 6.41
$ws.Range("E7").Value = 34.0
$ws.Range("F7").Value = 0.000000756

$ws.Range("C8").Value = 0.000597
$ws.Range("E8").Value = 49.0
$ws.Range("F8").Value = 0.0000000163

$ws.Range("C9").Value = 0.000278
$ws.Range("E9").Value = 22.8
$ws.Range("F9").Value = 0.0000231

$ws.Range("C10").Value = 0.000333
$ws.Range("E10").Value = 27.4
$ws.Range("F10").Value = 0.00000531

$ws.Range("C11").Value = 8.04
$ws.Range("E11").Value = 8.17
$ws.Range("F11").Value = 0.00667

$ws.Range("C12").Value = 2.67
$ws.Range("E12").Value = 2.71
$ws.Range("F12").Value = 0.107

$ws.Range("C13").Value = 0.926
$ws.Range("E13").Value = 0.94
$ws.Range("F13").Value = 0.338

$ws.Range("C14").Value = 1500.0
$ws.Range("E14").Value = 8.3
$ws.Range("F14").Value = 0.00629

$ws.Range("C15").Value = 8270.0
$ws.Range("E15").Value = 45.7
$ws.Range("F15").Value = 0.0000000357

$ws.Range("C16").Value = 667.0
$ws.Range("E16").Value = 3.69
$ws.Range("F16").Value = 0.0618

$ws.Range("C17").Value = 7.76
$ws.Range("E17").Value = 0.361
$ws.Range("F17").Value = 0.551

$ws.Range("C18").Value = 246.0
$ws.Range("E18").Value = 11.5
$ws.Range("F18").Value = 0.00157

$ws.Range("C19").Value = 26.7
$ws.Range("E19").Value = 1.24
$ws.Range("F19").Value = 0.271
